$d = $word.ActiveDocument

# The hidden "_GoBack" bookmark marks the last edit position. Remove it now;
# we'll re-add it at the right spot once the new text is in place so it
# doesn't get left behind at the wrong offset.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Anchor the edit to the "Version 1." text rather than a hard-coded offset.
$found = $d.Content
$found.Find.Execute("Version 1.", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($found.Start, $found.End)

# Rebuild that text as the run-split sequence produced by retyping mid-word
# ("Version" -> "Versi" | "on") and re-entering the version number and the
# trailing period ("1." -> "2" + "."), with the bookmark collapsed right
# before the re-typed ".".
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Versi</w:t></w:r><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)
